# Apply crypto price/volume updates per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.207.16'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.855.75'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7062'
$ws.Range('E5').Value = '  +2.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.06'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08027'
$ws.Range('E8').Value = '  +4.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3021'
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.45'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08181'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.854.60'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.185'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7023'
$ws.Range('E14').Value = '  -2.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.56'
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').Value = '29.121.72'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007959'
$ws.Range('E17').Value = '  +2.50%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.788'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.07'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = '2.049.64'
$ws.Range('E23').Value = '  -2.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.459'
$ws.Range('E24').Value = '  -1.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.77'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.898'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1429'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.07'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.921'
$ws.Range('E29').Value = '  -2.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.417'
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.476'
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.353'
$ws.Range('E32').Value = '  -2.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.026'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05184'
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.156'
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7136'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').Value = '  -2.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.644'
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01848'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.720'
$ws.Range('E40').Value = '  +1.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9370'
$ws.Range('E41').Value = '  +2.23%  '
$ws.Range('D42').Value = '1.127.75'
$ws.Range('E42').Value = '  +3.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.922'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4247'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.04'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.64'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5318'
$ws.Range('E48').Value = '  -4.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.758'
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.167'
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05953'
$ws.Range('E51').Value = '  +1.65%  '
